$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-9 and 11 per the regenerated save_data
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 10
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 2
$ws.Range("G11").Value = 1
